# Lesson_Names.xlsx - "Added 3eme-1ere I classes"
# Insert 8 new Original-Name/Abbreviation rows into the sorted B:C list on
# Tabelle1, keeping the existing alphabetical order intact, then move the
# active-cell selection the way the saved workbook shows it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows top-down using their FINAL row numbers: inserting a row
# at position N shifts everything currently at/after N down by one, which
# is exactly what's needed to open a gap at the next (larger) target row.
$newRows = @(
    @{ Row = 5; B = "Analyse et modélisation d'informations"; C = "Analyse" },
    @{ Row = 11; B = "Communication média"; C = "ComMédia" },
    @{ Row = 15; B = "Design graphique"; C = "Design" },
    @{ Row = 21; B = "Économie et finances"; C = "Economie" },
    @{ Row = 46; B = "Maîtrise d'ouvrage"; C = "Ouvrage" },
    @{ Row = 55; B = "Science de la programmation"; C = "Programmation" },
    @{ Row = 61; B = "Technologie et innovations"; C = "Technologie" },
    @{ Row = 62; B = "Technologies appliquées et projets"; C = "TechApp" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Rows("$r`:$r").Insert()
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
}

# Matches the saved selection recorded in the workbook after the edit.
$ws.Range("D14").Select()
